$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.855056881904602
$ws.Range("B1").Value = 3.061240196228027
$ws.Range("C1").Value = 2.974729537963867
$ws.Range("D1").Value = 3.429309368133545
$ws.Range("E1").Value = 3.970493793487549
